$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 2
$ws_ALC.Range("H2").Value = 740.6667

# ALC row 8
$ws_ALC.Range("H8").Value = 0
$ws_ALC.Range("I8").Value = 0
$ws_ALC.Range("K8").Value = 0
$ws_ALC.Range("M8").ClearContents()

# ALC row 15
$ws_ALC.Range("H15").Value = 310.2381
$ws_ALC.Range("I15").Value = 310.2381
$ws_ALC.Range("K15").Value = 930.7142999999999
$ws_ALC.Range("M15").Value = -761.7142999999999

# ALC row 19
$ws_ALC.Range("H19").Value = 1192.84
$ws_ALC.Range("I19").Value = 192.8
$ws_ALC.Range("J19").Value = 2692.9
$ws_ALC.Range("K19").Value = 192.8
$ws_ALC.Range("L19").Value = 2692.9
$ws_ALC.Range("M19").Value = -17.80000000000001
$ws_ALC.Range("N19").Value = -3042.9

# ALC row 33
$ws_ALC.Range("H33").Value = 156.72223
$ws_ALC.Range("I33").Value = 130.64706
$ws_ALC.Range("K33").Value = 130.64706
$ws_ALC.Range("M33").Value = 98.35293999999999

# ALC row 86
$ws_ALC.Range("H86").Value = 220600.8
$ws_ALC.Range("I86").Value = 1000
$ws_ALC.Range("J86").Value = 367001.34
$ws_ALC.Range("K86").Value = 1000
$ws_ALC.Range("L86").Value = 367001.34
$ws_ALC.Range("M86").Value = 123
$ws_ALC.Range("N86").Value = -369247.34

# ALC row 89
$ws_ALC.Range("H89").Value = 220600.8
$ws_ALC.Range("I89").Value = 1000
$ws_ALC.Range("J89").Value = 367001.34
$ws_ALC.Range("K89").Value = 5000
$ws_ALC.Range("L89").Value = 1835006.7
$ws_ALC.Range("M89").Value = 616
$ws_ALC.Range("N89").Value = -1846238.7

# ALC row 116
$ws_ALC.Range("H116").Value = 3013.4285
$ws_ALC.Range("I116").Value = 3019
$ws_ALC.Range("J116").Value = 2999.5
$ws_ALC.Range("K116").Value = 3019
$ws_ALC.Range("L116").Value = 2999.5
$ws_ALC.Range("M116").Value = 423
$ws_ALC.Range("N116").Value = -9883.5

# ALC row 137
$ws_ALC.Range("H137").Value = 1944.6154
$ws_ALC.Range("I137").Value = 1186.7778
$ws_ALC.Range("K137").Value = 3560.3334
$ws_ALC.Range("M137").Value = -1010.3334

# ARM row 4
$ws_ARM.Range("H4").Value = 298
$ws_ARM.Range("I4").Value = 298
$ws_ARM.Range("K4").Value = 298
$ws_ARM.Range("M4").Value = -182

# ARM row 5
$ws_ARM.Range("H5").Value = 995
$ws_ARM.Range("I5").Value = 975
$ws_ARM.Range("K5").Value = 975
$ws_ARM.Range("M5").Value = -863

# ARM row 74
$ws_ARM.Range("H74").Value = 2262.5454
$ws_ARM.Range("I74").Value = 1612
$ws_ARM.Range("K74").Value = 1612
$ws_ARM.Range("M74").Value = -738

# ARM row 77
$ws_ARM.Range("H77").Value = 2262.5454
$ws_ARM.Range("I77").Value = 1612
$ws_ARM.Range("K77").Value = 8060
$ws_ARM.Range("M77").Value = -3692

# ARM row 132
$ws_ARM.Range("H132").Value = 6396.727
$ws_ARM.Range("I132").Value = 5395.3335
$ws_ARM.Range("J132").Value = 7598.4
$ws_ARM.Range("K132").Value = 16186.0005
$ws_ARM.Range("L132").Value = 22795.2
$ws_ARM.Range("M132").Value = -13656.0005
$ws_ARM.Range("N132").Value = -27855.2

# BSM row 4
$ws_BSM.Range("H4").Value = 995
$ws_BSM.Range("I4").Value = 975
$ws_BSM.Range("K4").Value = 975
$ws_BSM.Range("M4").Value = -860

# BSM row 86
$ws_BSM.Range("H86").Value = 1166.3334
$ws_BSM.Range("I86").Value = 999
$ws_BSM.Range("K86").Value = 999
$ws_BSM.Range("M86").Value = 124

# BSM row 89
$ws_BSM.Range("H89").Value = 1166.3334
$ws_BSM.Range("I89").Value = 999
$ws_BSM.Range("K89").Value = 4995
$ws_BSM.Range("M89").Value = 621

# BSM row 134
$ws_BSM.Range("H134").Value = 3027.3076
$ws_BSM.Range("I134").Value = 2350.4546
$ws_BSM.Range("K134").Value = 7051.3638
$ws_BSM.Range("M134").Value = -4516.3638

# CRP row 31
$ws_CRP.Range("H31").Value = 2684.8147
$ws_CRP.Range("I31").Value = 1643.5714
$ws_CRP.Range("J31").Value = 6329.1665
$ws_CRP.Range("K31").Value = 1643.5714
$ws_CRP.Range("L31").Value = 6329.1665
$ws_CRP.Range("M31").Value = -1348.5714
$ws_CRP.Range("N31").Value = -6919.1665

# CRP row 34
$ws_CRP.Range("H34").Value = 2684.8147
$ws_CRP.Range("I34").Value = 1643.5714
$ws_CRP.Range("J34").Value = 6329.1665
$ws_CRP.Range("K34").Value = 1643.5714
$ws_CRP.Range("L34").Value = 6329.1665
$ws_CRP.Range("M34").Value = -1441.5714
$ws_CRP.Range("N34").Value = -6733.1665

# CRP row 58
$ws_CRP.Range("H58").Value = 3411.25
$ws_CRP.Range("I58").Value = 1949
$ws_CRP.Range("K58").Value = 1949
$ws_CRP.Range("M58").Value = -1746

# CRP row 132
$ws_CRP.Range("H132").Value = 2010.0625
$ws_CRP.Range("I132").Value = 1904.3572
$ws_CRP.Range("K132").Value = 5713.071599999999
$ws_CRP.Range("M132").Value = -3183.071599999999

# CRP row 136
$ws_CRP.Range("H136").Value = 3411.25
$ws_CRP.Range("I136").Value = 1949
$ws_CRP.Range("K136").Value = 5847
$ws_CRP.Range("M136").Value = -3297

# CUL row 25
$ws_CUL.Range("H25").Value = 225
$ws_CUL.Range("I25").Value = 216.66667
$ws_CUL.Range("K25").Value = 650.00001
$ws_CUL.Range("M25").Value = -481.00001

# CUL row 30
$ws_CUL.Range("H30").Value = 225
$ws_CUL.Range("I30").Value = 216.66667
$ws_CUL.Range("K30").Value = 650.00001
$ws_CUL.Range("M30").Value = -548.00001

# CUL row 34
$ws_CUL.Range("H34").Value = 1954.1428
$ws_CUL.Range("J34").Value = 3375
$ws_CUL.Range("L34").Value = 10125
$ws_CUL.Range("N34").Value = -10293

# CUL row 39
$ws_CUL.Range("H39").Value = 11569.429
$ws_CUL.Range("I39").Value = 0
$ws_CUL.Range("J39").Value = 11569.429
$ws_CUL.Range("K39").Value = 0
$ws_CUL.Range("L39").Value = 34708.287
$ws_CUL.Range("M39").ClearContents()
$ws_CUL.Range("N39").Value = -35296.287

# CUL row 55
$ws_CUL.Range("H55").Value = 7372.5454
$ws_CUL.Range("J55").Value = 9174.875
$ws_CUL.Range("L55").Value = 27524.625
$ws_CUL.Range("N55").Value = -27878.625

# CUL row 131
$ws_CUL.Range("H131").Value = 1376.409
$ws_CUL.Range("I131").Value = 1208
$ws_CUL.Range("K131").Value = 3624
$ws_CUL.Range("M131").Value = 1416

# GSM row 102
$ws_GSM.Range("H102").Value = 1292.579
$ws_GSM.Range("I102").Value = 1304.2142
$ws_GSM.Range("K102").Value = 1304.2142
$ws_GSM.Range("M102").Value = 317.7858000000001

# GSM row 122
$ws_GSM.Range("H122").Value = 70493.87
$ws_GSM.Range("I122").Value = 3097.1
$ws_GSM.Range("K122").Value = 9291.299999999999
$ws_GSM.Range("M122").Value = -6841.299999999999

# GSM row 132
$ws_GSM.Range("H132").Value = 3975.4092
$ws_GSM.Range("I132").Value = 3975.4092
$ws_GSM.Range("J132").Value = 0
$ws_GSM.Range("K132").Value = 11926.2276
$ws_GSM.Range("L132").Value = 0
$ws_GSM.Range("M132").Value = -9396.2276
$ws_GSM.Range("N132").ClearContents()

# LTW row 46
$ws_LTW.Range("H46").Value = 1884.875
$ws_LTW.Range("I46").Value = 834.75
$ws_LTW.Range("J46").Value = 2935
$ws_LTW.Range("K46").Value = 834.75
$ws_LTW.Range("L46").Value = 2935
$ws_LTW.Range("M46").Value = -646.75
$ws_LTW.Range("N46").Value = -3311

# LTW row 82
$ws_LTW.Range("H82").Value = 84524.25
$ws_LTW.Range("I82").Value = 961.75
$ws_LTW.Range("K82").Value = 961.75
$ws_LTW.Range("M82").Value = -600.75

# LTW row 85
$ws_LTW.Range("H85").Value = 84524.25
$ws_LTW.Range("I85").Value = 961.75
$ws_LTW.Range("K85").Value = 961.75
$ws_LTW.Range("M85").Value = 286.25

# LTW row 93
$ws_LTW.Range("H93").Value = 738.5714
$ws_LTW.Range("I93").Value = 738.5714
$ws_LTW.Range("K93").Value = 738.5714
$ws_LTW.Range("M93").Value = 509.4286

# LTW row 132
$ws_LTW.Range("H132").Value = 85917.086
$ws_LTW.Range("I132").Value = 102400.5
$ws_LTW.Range("K132").Value = 307201.5
$ws_LTW.Range("M132").Value = -304671.5

# LTW row 136
$ws_LTW.Range("H136").Value = 7517.6
$ws_LTW.Range("I136").Value = 7216
$ws_LTW.Range("K136").Value = 21648
$ws_LTW.Range("M136").Value = -19098

# WVR row 62
$ws_WVR.Range("H62").Value = 0
$ws_WVR.Range("I62").Value = 0
$ws_WVR.Range("K62").Value = 0
$ws_WVR.Range("M62").ClearContents()

# WVR row 65
$ws_WVR.Range("H65").Value = 0
$ws_WVR.Range("I65").Value = 0
$ws_WVR.Range("K65").Value = 0
$ws_WVR.Range("M65").ClearContents()

# WVR row 136
$ws_WVR.Range("H136").Value = 6171
$ws_WVR.Range("I136").Value = 6057
$ws_WVR.Range("K136").Value = 18171
$ws_WVR.Range("M136").Value = -15621
